# Slide 9 (1-based) / Shape 1 is the title placeholder
# "Google Shape;236;p37" that currently reads "Industrial Revolution (IR2)".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(1)

# Widen the title text box slightly to fit the longer text
# (5625610 EMU -> 6039200 EMU, i.e. 442.961... pt -> 475.5276 pt).
$sh.Width = 475.5276

# Re-word "IR2" into "IR 2.0" by replacing just the trailing
# "IR2)" characters; PowerPoint splits this into its own run.
$tr = $sh.TextFrame.TextRange
$tail = $tr.Characters(24, 4)
$tail.Text = "IR 2.0)"
